# Link Triggered BP.xlsx - data refresh for "Feb" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feb")

# Fix the header typo "Divfision" -> "Division" (now matches other sheets)
$ws.Range("F1").Value = "Division"

# Updated counts / figures
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 45

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 34.48

$ws.Range("B7").Value = 39
$ws.Range("E7").Value = 51.28

$ws.Range("B9").Value = 17
$ws.Range("B10").Value = 14
$ws.Range("B11").Value = 19
$ws.Range("B13").Value = 10

$ws.Range("B14").Value = 40
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 25

$ws.Range("B15").Value = 42
$ws.Range("E15").Value = 23.81

# Update the active selection to match the refreshed view
$ws.Range("F10:F11").Select() | Out-Null
